# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary rows (9-12): update row-label styles + score numbers ---
# A10/A11/A12 switch from default (Normal) style to the "mtitleStyle" style,
# matching the look of the other header rows (e.g. row 9).
$ws.Range("A9").Copy()
foreach ($r in @(10,11,12)) {
  $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
}

# Row 10 "No." counts
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 28

# Row 11 "Marking" weights (C11 becomes a real number instead of text "-1")
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 "Total" row
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "66/112"

# --- Answer grid (rows 16-40) ---
# Column A becomes the "Student Ans" column for the first question block,
# mirroring column B ("Correct Ans"): green/correctStyle when it matches,
# red/incorrectStyle when it doesn't, and left blank (normalStyle) when the
# student did not attempt the question.

# Reusable format sources already present in the sheet:
#   B10 -> correctStyle (green)
#   C10 -> incorrectStyle (red)
$correctCells = @(16,18,19,22,23,25,26,27,28,30,31,33,36,39,40)
$correctAnswers = @{
  16 = "Option A"; 18 = "Option B"; 19 = "Option C"; 22 = "Option D";
  23 = "Option D"; 25 = "Option A"; 26 = "Option C"; 27 = "Option A";
  28 = "Option D"; 30 = "Option B"; 31 = "Option D"; 33 = "Option D";
  36 = "Option A"; 39 = "Option D"; 40 = "Option D"
}

$ws.Range("B10").Copy()
foreach ($r in $correctCells) {
  $cell = $ws.Range("A$r")
  $cell.PasteSpecial(-4122) | Out-Null
  $cell.Value = $correctAnswers[$r]
}

$wrongAnswers = @{ 32 = "Option B"; 38 = "Option C" }
$ws.Range("C10").Copy()
foreach ($r in @(32,38)) {
  $cell = $ws.Range("A$r")
  $cell.PasteSpecial(-4122) | Out-Null
  $cell.Value = $wrongAnswers[$r]
}

# Second question block (columns D/E) is now only 2 questions long; row 16/17
# get real "Student Ans" values in column D, formatted the same as column A.
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Value = "Option C"

# The rest of the D/E block (rows 19-40) and the entire third question block
# (columns G/H, rows 15-21) belonged to question sets that no longer exist.
$ws.Range("D19:E40").Clear()
$ws.Columns("G:H").Delete()

$wb.Save()
